$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D37").Value = "Importante"
$ws.Range("E37").Value = "Funcional"
$ws.Range("F37").Value = "Media"
$ws.Range("G37").Value = "v2"

$ws.Range("D38").Value = "Importante"
$ws.Range("E38").Value = "Funcional"
$ws.Range("F38").Value = "Media"
$ws.Range("G38").Value = "v2"

$ws.Range("D39").Value = "Importante"
$ws.Range("E39").Value = "Funcional"
$ws.Range("F39").Value = "Media"
$ws.Range("G39").Value = "v2"

$ws.Range("D40").Value = "Opcional"
$ws.Range("E40").Value = "Funcional"
$ws.Range("F40").Value = "Difícil"
$ws.Range("G40").Value = "v3"

$ws.Range("D41").Value = "Mínimo"
$ws.Range("E41").Value = "Funcional"
$ws.Range("F41").Value = "Media"
$ws.Range("G41").Value = "v3"

$ws.Range("D42").Value = "Mínimo"
$ws.Range("E42").Value = "Funcional"
$ws.Range("F42").Value = "Media"
$ws.Range("G42").Value = "v3"

$ws.Range("D43").Value = "Importante"
$ws.Range("E43").Value = "Funcional"
$ws.Range("F43").Value = "Fácil"
$ws.Range("G43").Value = "v2"

$ws.Range("D44").Value = "Mínimo"
$ws.Range("E44").Value = "Funcional"
$ws.Range("F44").Value = "Media"
$ws.Range("G44").Value = "v3"

$ws.Range("D45").Value = "Mínimo"
$ws.Range("E45").Value = "Funcional"
$ws.Range("F45").Value = "Fácil"
$ws.Range("G45").Value = "v2"

$ws.Range("D46").Value = "Importante"
$ws.Range("E46").Value = "Funcional"
$ws.Range("F46").Value = "Fácil"
$ws.Range("G46").Value = "v3"

$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Range("F46").Select()
